$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target: cell B11 changes from the text "R40" to the text "1" (a
# shared string), keeping its existing cell style (s="23") untouched.
#
# A plain `Range.Value = "1"` would be auto-typed as a NUMBER by Excel
# (since "1" parses as numeric), which also changes the cell's applied
# number format / style index - not what we want.
#
# Instead, build the literal text "1" as the result of a formula
# (="1") in a scratch cell outside the used range - formula results
# keep their string type regardless of how "number-like" they look -
# then copy just the computed value over to B11 and clean up the
# scratch cell so it leaves no trace in the saved workbook.
$helper = $ws.Range("ZZ1")
$helper.Formula = "=""1"""
$helper.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()
